$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 176, shifting existing rows 176-280 down to 177-281
$ws.Rows.Item(176).Insert()

# Populate the newly inserted row 176 with its data
$ws.Range("A176").Value = 3
$ws.Range("B176").Value = "Femacal de La Calera"
$ws.Range("C176").Value = "Coquimbo"
$ws.Range("D176").Value = "2022-01-21"
$ws.Range("E176").Value = 5
$ws.Range("F176").Value = 100114013
$ws.Range("G176").Value = "Zanahoria"
$ws.Range("H176").Value = "Sin especificar"
$ws.Range("I176").Value = "Primera"
$ws.Range("J176").Value = 760
$ws.Range("K176").Value = 7000
$ws.Range("L176").Value = 7500
$ws.Range("M176").Value = 7250
$ws.Range("N176").Value = "$/saco 20 kilos"
$ws.Range("O176").Value = "Provincia de Quillota"
$ws.Range("P176").Value = 362
$ws.Range("Q176").Value = 20
$ws.Range("R176").Value = "Hortaliza"
